$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original row 3 values (Alan Moreno / password123 / 3) before shifting
$origA3 = $ws.Range("A3").Value2
$origB3 = $ws.Range("B3").Value2
$origC3 = $ws.Range("C3").Value2

# Shift rows 4,5,6 up into rows 3,4,5 (per-cell, since Value returns a non-expandable object for multi-cell ranges)
for ($col = 1; $col -le 3; $col++) {
    $ws.Cells.Item(3, $col).Value = $ws.Cells.Item(4, $col).Value2
}
for ($col = 1; $col -le 3; $col++) {
    $ws.Cells.Item(4, $col).Value = $ws.Cells.Item(5, $col).Value2
}
for ($col = 1; $col -le 3; $col++) {
    $ws.Cells.Item(5, $col).Value = $ws.Cells.Item(6, $col).Value2
}

# Place original row3 name/password at row6 with updated access level
$ws.Range("A6").Value = $origA3
$ws.Range("B6").Value = $origB3
$ws.Range("C6").Value = 2

# Update the active selection to F3, matching the saved view state
$null = $ws.Range("F3").Select()
